$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EJ44")

# New header labels for the "Extra Hard" EJ44 string set
$ws.Range("B1").Value = "J4401"
$ws.Range("C1").Value = "J4402"
$ws.Range("D1").Value = "J4403"
$ws.Range("E1").Value = "J4404"
$ws.Range("F1").Value = "J4405"
$ws.Range("G1").Value = "J4406"

# New tension figures for EJ44 (Extra Hard)
$ws.Range("B2").Value = 293.8
$ws.Range("C2").Value = 220.1
$ws.Range("D2").Value = 177.7
$ws.Range("E2").Value = 132.7
$ws.Range("F2").Value = 97.7
$ws.Range("G2").Value = 74.5

$ws.Range("B3").Value = 299.8
$ws.Range("C3").Value = 227.3
$ws.Range("D3").Value = 182.5
$ws.Range("E3").Value = 135.2
$ws.Range("F3").Value = 99.8
$ws.Range("G3").Value = 76.3

$ws.Range("B4").Value = 307.5
$ws.Range("C4").Value = 235.8
$ws.Range("D4").Value = 188.1
$ws.Range("E4").Value = 138.7
$ws.Range("F4").Value = 101.7
$ws.Range("G4").Value = 78.1

$ws.Range("B5").Value = 315.2
$ws.Range("C5").Value = 242.7
$ws.Range("D5").Value = 192.8
$ws.Range("E5").Value = 142
$ws.Range("F5").Value = 104.1
$ws.Range("G5").Value = 80

$ws.Range("B6").Value = 322
$ws.Range("C6").Value = 248.5
$ws.Range("D6").Value = 196.9
$ws.Range("E6").Value = 144.8
$ws.Range("F6").Value = 106.1
$ws.Range("G6").Value = 81.7

$ws.Range("B7").Value = 328.3
$ws.Range("C7").Value = 254.1
$ws.Range("D7").Value = 201.2
$ws.Range("E7").Value = 147.3
$ws.Range("F7").Value = 108.2
$ws.Range("G7").Value = 83.5

$ws.Range("B8").Value = 334.7
$ws.Range("C8").Value = 260.3
$ws.Range("D8").Value = 206.4
$ws.Range("E8").Value = 150.5
$ws.Range("F8").Value = 110.5
$ws.Range("G8").Value = 85.2

# Make EJ44 the active sheet / tab, with the selection left where Excel
# would land after entering the last value (G9), matching the saved view.
$ws.Activate()
$ws.Range("G9").Select()
